# Daily attendance processing - 2026-01-09 22:35:03
# Swap the order of "Recorded By" names in column G from
# "System, dnasr281@gmail.com" to "dnasr281@gmail.com, System"
# for the specific rows whose value currently reads "System, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetCells = @(
    "G2","G3","G4","G5","G6","G7",
    "G16","G17",
    "G22","G23",
    "G37","G38",
    "G43","G44",
    "G58","G59",
    "G64","G65",
    "G79","G80",
    "G85","G86","G87","G88","G89","G90",
    "G99","G100",
    "G105","G106","G107","G108","G109","G110",
    "G119","G120",
    "G125","G126","G127","G128","G129","G130",
    "G139","G140",
    "G145","G146","G147","G148","G149","G150",
    "G159","G160",
    "G165","G166","G167","G168","G169","G170",
    "G179","G180",
    "G185","G186",
    "G200","G201",
    "G206","G207",
    "G221","G222",
    "G227","G228",
    "G242","G243"
)

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

foreach ($addr in $targetCells) {
    $cell = $ws.Range($addr)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
